# tix-6 added read balance method
#
# Adds a new "Checkings" charge transaction (store: "rich peoiple store",
# charge: 9999, new balance: 89978.2) dated 2021/10/01 20:05:14 to both the
# Home overview sheet and the Checking account sheet, and updates the
# Checking account running balance on the Data sheet accordingly.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Writing a date/time-shaped literal string via .Value causes Excel to
    # auto-convert it into a date/time serial number. Route the literal
    # through a TEXT() formula and then paste-special just the resulting
    # value back over the formula so the cell ends up holding a plain
    # shared-string, exactly like the other text cells in this workbook.
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

# ---------------------------------------------------------------------
# Home sheet: append the new transaction as row 30
# ---------------------------------------------------------------------
$home = $wb.Worksheets.Item("Home")
Set-TextValue $home.Range("F30") "2021/10/01"
Set-TextValue $home.Range("G30") "20:05:14"
$home.Range("H30").Value = "Checkings"
$home.Range("I30").Value = 9999
Set-TextValue $home.Range("J30") "rich peoiple store"
$home.Range("K30").Value = 89978.2

# ---------------------------------------------------------------------
# Checking sheet: append the same transaction as row 16
# ---------------------------------------------------------------------
$checking = $wb.Worksheets.Item("Checking")
Set-TextValue $checking.Range("F16") "2021/10/01"
Set-TextValue $checking.Range("G16") "20:05:14"
$checking.Range("H16").Value = "Checkings"
$checking.Range("I16").Value = 9999
Set-TextValue $checking.Range("J16") "rich peoiple store"
$checking.Range("K16").Value = 89978.2

# ---------------------------------------------------------------------
# Data sheet: roll the Checking account running balance forward
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")
$data.Range("A3").Value = 89978.2
